$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.360dx.com/cancer/epredia-gains-fda-510k-digital-pathology-imaging-system"
$keyword = "digital pathology"
$title = "Epredia Gains FDA 510(k) for Digital Pathology Imaging System"

$ws.Range("A23").Value = $url
$ws.Range("B23").Value = $keyword
$ws.Range("C23").Value = $title

$ws.Hyperlinks.Add($ws.Range("A23"), $url)
$ws.Range("A23").Style = $ws.Range("A22").Style
